# Daily attendance processing - 2025-10-25 19:41:37
# Reorders the "Recorded By" (column G) entries on the "Session Analysis
# Results" sheet: any value beginning with the literal prefix "System, "
# has that leading "System, " removed and ", System" appended instead,
# effectively moving "System" from the front of the list to the back.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$prefix = "System, "
$used = $ws.UsedRange
$firstRow = $used.Row
$lastRow = $firstRow + $used.Rows.Count - 1

$changed = 0

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $current = $cell.Value2

    if ($current -ne $null -and $current.GetType().Name -eq "String" -and $current.StartsWith($prefix)) {
        $rest = $current.Substring($prefix.Length)
        $newValue = $rest + ", System"
        $cell.Value = $newValue
        $changed++
    }
}

Write-Output "Rows updated: $changed"
